# Add "Mandatory" column to the structures sheets (Clinical Structures,
# opt structures, couch_structures) and mark most rows as "oui" (mandatory),
# plus add a new "manque" row to "opt structures". Also update the
# active-sheet/selection bookkeeping to match the saved workbook state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# General data (sheet1): loses the active-tab/selection it used to have.
# Do this first so that a later Activate() on another sheet "wins" and
# ends up as the final active tab (Range.Select implicitly activates its
# own worksheet, so ordering matters).
# ---------------------------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("General data")
$wsGeneral.Range("C21").Select()

# ---------------------------------------------------------------------
# Clinical Structures (sheet2): column G header + "oui" markers
# ---------------------------------------------------------------------
$wsClinical = $wb.Worksheets.Item("Clinical Structures")
$wsClinical.Range("G1").Value = "Mandatory"
$wsClinical.Range("G10").Value = "oui"
$wsClinical.Range("G15").Value = "oui"
$wsClinical.Range("G17").Value = "oui"
$wsClinical.Range("G19").Value = "oui"
$wsClinical.Range("G20").Value = "oui"
$wsClinical.Range("G21").Value = "oui"
$wsClinical.Range("G23").Value = "oui"
$wsClinical.Range("G24").Select()

# ---------------------------------------------------------------------
# opt structures (sheet3): column G header, "oui" markers, and a new
# "manque" row at the bottom
# ---------------------------------------------------------------------
$wsOpt = $wb.Worksheets.Item("opt structures")
$wsOpt.Range("G1").Value = "Mandatory"
$wsOpt.Range("G6").Value = "oui"
$wsOpt.Range("G7").Value = "oui"
$wsOpt.Range("G23").Value = "oui"
$wsOpt.Range("G28").Value = "oui"
$wsOpt.Range("A29").Value = "manque"
$wsOpt.Range("G29").Value = "oui"
$wsOpt.Range("A29").Select()

# ---------------------------------------------------------------------
# couch_structures (sheet4): column G header + "oui" markers
# ---------------------------------------------------------------------
$wsCouch = $wb.Worksheets.Item("couch_structures")
$wsCouch.Range("G1").Value = "Mandatory"
$wsCouch.Range("G2").Value = "oui"
$wsCouch.Range("G3").Value = "oui"
$wsCouch.Range("G4").Value = "oui"
$wsCouch.Range("G5").Value = "oui"
$wsCouch.Range("G6").Select()

# This sheet becomes the active tab (activeTab="3", tabSelected="1").
# Activate() (and the Select() just above, on this same sheet) must be
# the last selection-affecting calls so this ends up as the final active
# sheet/selection in the saved workbook.
$wsCouch.Activate()
